$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "effects" column from the TREATMENT table.
# Before: D27=effects, E27=startDate, F27=endDate (and F28:F31 blank placeholders)
# After:  D27=startDate, E27=endDate (F27:F31 removed entirely)
$ws.Range("D27").Value = "startDate"
$ws.Range("E27").Value = "endDate"
$ws.Range("F27:F31").Clear()

# Update the view: scroll so row 18 is at the top, and select G30.
$ws.Range("G30").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 18
